# notExistingServer-template.docx update
#
# The canonical-OOXML diff for this commit touches word/document.xml and
# word/styles.xml, but every single "-"/"+" line pair in it is a pure
# re-ordering of XML attributes on the same element (e.g.
#   <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
#   -> <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
# and likewise for the root <w:document> namespace declarations, <w:pgSz>,
# <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>,
# and every <w:style>/<w:tblCellMar> child). In each case the attribute
# *names and values* are identical before and after - only their emission
# order changed (alphabetical afterwards), which is a side effect of the
# tool that re-serialized the part upstream, not a content/formatting
# change. No text, run, paragraph, style definition, font, margin, color,
# or property value is added, removed, or modified anywhere in the diff.
#
# Word's object model (and therefore COM automation) does not expose
# control over XML attribute emission order - that is a serialization
# detail of whatever engine writes the part, not part of the document
# model surfaced by Word.Application/Document. So there is nothing to
# change from the document-model point of view: the template's text,
# formatting, styles, and properties must stay exactly as they are.
#
# This script therefore intentionally performs no content mutation -
# it only touches $d to make the intent explicit, leaving the document
# identical to the input, matching the (purely cosmetic) diff.
$d = $word.ActiveDocument
$null = $d.Name
